# Update the multiplication problems in the practice-sheet table.
# Each cell contains a run of text like "49×67=" that is replaced with a
# new problem, e.g. "27×12=". All values are unique within the document,
# so a simple Find/Replace for each exact pair is safe and unambiguous.

$d = $word.ActiveDocument

$replacements = @(
    @("49×67=", "27×12="),
    @("89×86=", "15×38="),
    @("35×86=", "21×20="),
    @("36×45=", "37×27="),
    @("34×25=", "38×59="),
    @("89×79=", "62×64="),
    @("53×30=", "11×30="),
    @("83×73=", "31×75="),
    @("17×21=", "65×89="),
    @("21×80=", "78×33="),
    @("76×39=", "32×32="),
    @("65×25=", "48×53="),
    @("20×56=", "61×39="),
    @("96×38=", "87×37="),
    @("88×12=", "20×25="),
    @("97×35=", "36×20="),
    @("31×40=", "41×38="),
    @("64×99=", "38×21="),
    @("60×79=", "45×40="),
    @("95×86=", "66×40="),
    @("91×16=", "89×76="),
    @("54×66=", "14×57="),
    @("83×70=", "87×21="),
    @("47×74=", "39×95="),
    @("91×38=", "63×49=")
)

foreach ($pair in $replacements) {
    $oldText = $pair[0]
    $newText = $pair[1]

    $found = $d.Content.Find.Execute(
        $oldText, $true, $false, $false, $false, $false,
        $true, 1, $false, $newText, 2)

    if (-not $found) {
        Write-Host "WARNING: could not find '$oldText' to replace with '$newText'"
    }
}

Write-Host "Done replacing $($replacements.Count) multiplication problems."
